$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host ("Sheet: " + $ws.Name)
Write-Host ("A1: " + $ws.Range("A1").Value())
Write-Host ("H53: " + $ws.Range("H53").Value())
Write-Host ("C54: " + $ws.Range("C54").Value())
Write-Host ("H54: " + $ws.Range("H54").Value())
Write-Host ("F3: " + $ws.Range("F3").Value())
Write-Host ("F43: " + $ws.Range("F43").Value())
Write-Host ("F54: " + $ws.Range("F54").Value())
